$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds CorpCode values stored as text. Force text format so the
# new numeric-looking strings stay text (matching the original inlineStr
# string cells) rather than being coerced into numbers.
$range = $ws.Range("A2:A23")
$range.NumberFormat = "@"

$ws.Range("A2").Value = "1543425"
$ws.Range("A3").Value = "1546848"
$ws.Range("A4").Value = "1550270"
$ws.Range("A5").Value = "1553693"
$ws.Range("A6").Value = "1557116"
$ws.Range("A7").Value = "1560538"
$ws.Range("A8").Value = "1567380"
$ws.Range("A9").Value = "1591335"
$ws.Range("A10").Value = "1601602"
$ws.Range("A11").Value = "1618714"
$ws.Range("A12").Value = "1639249"
$ws.Range("A13").Value = "1642669"
$ws.Range("A14").Value = "1646092"
$ws.Range("A15").Value = "1656359"
$ws.Range("A16").Value = "1659782"
$ws.Range("A17").Value = "1670046"
$ws.Range("A18").Value = "1673469"
$ws.Range("A19").Value = "1680315"
$ws.Range("A20").Value = "1683738"
$ws.Range("A21").Value = "1690583"
$ws.Range("A22").Value = "1694002"
$ws.Range("A23").Value = "1707692"
